# Applies the edits described in the commit diff:
#   - Inserts a new match row (Azerbaijan Premier League) at row 4,
#     pushing the existing rows 4-13 down to rows 5-14.
#   - Updates several odds/time values across the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new row at position 4 (shifts existing rows 4-13 down to 5-14)
$ws.Rows.Item(4).Insert()

# 2) Populate the newly inserted row 4 with the Azerbaijan Premier League match data.
#    The "Date" value (column B) is copied from an existing date cell (rather than
#    assigned as a literal string) so Excel keeps it stored as text instead of
#    auto-converting the "yyyy-mm-dd" text into a serial date value.
$ws.Range("B5").Copy()
$ws.Range("B4").PasteSpecial(-4163)
$ws.Cells.Item(4, 1).Value = "Azerbaijan Premier League"
$ws.Cells.Item(4, 3).Value = "07:30:00"
$ws.Cells.Item(4, 4).Value = "Karvan Evlakh"
$ws.Cells.Item(4, 5).Value = "FK Sumqayit"
$ws.Cells.Item(4, 6).Value = 1.04
$ws.Cells.Item(4, 7).Value = 1000
$ws.Cells.Item(4, 8).Value = 1.04
$ws.Cells.Item(4, 9).Value = 1000
$ws.Cells.Item(4, 10).Value = 1.02
$ws.Cells.Item(4, 11).Value = 950
$ws.Cells.Item(4, 12).Value = 1.01
$ws.Cells.Item(4, 13).Value = 1.01
$ws.Cells.Item(4, 14).Value = 1.24
$ws.Cells.Item(4, 15).Value = 1.01
$ws.Cells.Item(4, 16).Value = 1.24
$ws.Cells.Item(4, 17).Value = 1.01
$ws.Cells.Item(4, 18).Value = 1.07
$ws.Cells.Item(4, 19).Value = 1.01
$ws.Cells.Item(4, 20).Value = 1.01
$ws.Cells.Item(4, 21).Value = 1.01
$ws.Cells.Item(4, 22).Value = 1.01
$ws.Cells.Item(4, 23).Value = 1.01
$ws.Cells.Item(4, 24).Value = 1000
$ws.Cells.Item(4, 25).Value = 1000
$ws.Cells.Item(4, 26).Value = 1000
$ws.Cells.Item(4, 27).Value = 1000
$ws.Cells.Item(4, 28).Value = 1000
$ws.Cells.Item(4, 29).Value = 1000
$ws.Cells.Item(4, 30).Value = 1000
$ws.Cells.Item(4, 31).Value = 1000
$ws.Cells.Item(4, 32).Value = 1000
$ws.Cells.Item(4, 33).Value = 1000
$ws.Cells.Item(4, 34).Value = 1000
$ws.Cells.Item(4, 35).Value = 1000
$ws.Cells.Item(4, 36).Value = 1000
$ws.Cells.Item(4, 37).Value = 1000
$ws.Cells.Item(4, 38).Value = 1000
$ws.Cells.Item(4, 39).Value = 1000
$ws.Cells.Item(4, 40).Value = 1000
$ws.Cells.Item(4, 41).Value = 1000

# 3) Apply value corrections to row 2 (unaffected by the row insertion)
$ws.Cells.Item(2, 16).Value = 1.98

# 4) Apply value corrections to row 3 (unaffected by the row insertion)
$ws.Cells.Item(3, 3).Value = "06:05:00"
$ws.Cells.Item(3, 6).Value = 2.54
$ws.Cells.Item(3, 7).Value = 2.66
$ws.Cells.Item(3, 8).Value = 2.66
$ws.Cells.Item(3, 9).Value = 2.78
$ws.Cells.Item(3, 13).Value = 1.03
$ws.Cells.Item(3, 18).Value = 1.65
$ws.Cells.Item(3, 22).Value = 1.56
$ws.Cells.Item(3, 23).Value = 1.6
$ws.Cells.Item(3, 38).Value = 980

# 5) Apply value corrections to the rows that were shifted down by the insertion
# Row 6 (originally row 5)
$ws.Cells.Item(6, 8).Value = 1.92
# Row 7 (originally row 6)
$ws.Cells.Item(7, 10).Value = 3.75
# Row 9 (originally row 8)
$ws.Cells.Item(9, 6).Value = 1.94
$ws.Cells.Item(9, 8).Value = 4.8
$ws.Cells.Item(9, 9).Value = 7.4
$ws.Cells.Item(9, 10).Value = 2.6
$ws.Cells.Item(9, 11).Value = 3.65
$ws.Cells.Item(9, 17).Value = 3.15
# Row 10 (originally row 9)
$ws.Cells.Item(10, 6).Value = 3
$ws.Cells.Item(10, 7).Value = 3.25
$ws.Cells.Item(10, 8).Value = 2.8
$ws.Cells.Item(10, 9).Value = 2.98
$ws.Cells.Item(10, 10).Value = 2.88
# Row 12 (originally row 11)
$ws.Cells.Item(12, 14).Value = 4.8
$ws.Cells.Item(12, 28).Value = 9.199999999999999
# Row 13 (originally row 12)
$ws.Cells.Item(13, 14).Value = 2.62
# Row 14 (originally row 13)
$ws.Cells.Item(14, 14).Value = 2.72
